# Auto-generated edit script: updates market-price derived columns (H-N)
# for the rows identified in the commit diff, across all 8 item-category sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4666.6665
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 6250
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 6250
$ws.Range("M32").Value = -1174
$ws.Range("N32").Value = -6902
$ws.Range("H38").Value = 2571.2727
$ws.Range("I38").Value = 941.1111
$ws.Range("J38").Value = 9907
$ws.Range("K38").Value = 2823.3333
$ws.Range("L38").Value = 29721
$ws.Range("M38").Value = -2451.3333
$ws.Range("N38").Value = -30465
$ws.Range("H40").Value = 9850.4
$ws.Range("J40").Value = 14500
$ws.Range("L40").Value = 14500
$ws.Range("N40").Value = -14850
$ws.Range("H86").Value = 7521668.5
$ws.Range("J86").Value = 10527975
$ws.Range("L86").Value = 10527975
$ws.Range("N86").Value = -10530221
$ws.Range("H89").Value = 7521668.5
$ws.Range("J89").Value = 10527975
$ws.Range("L89").Value = 52639875
$ws.Range("N89").Value = -52651107
$ws.Range("H92").Value = 157.09091
$ws.Range("I92").Value = 125.44444
$ws.Range("K92").Value = 125.44444
$ws.Range("M92").Value = 1122.55556
$ws.Range("H107").Value = 34543.6
$ws.Range("I107").Value = 44339.523
$ws.Range("J107").Value = 2357
$ws.Range("K107").Value = 44339.523
$ws.Range("L107").Value = 2357
$ws.Range("M107").Value = -42419.523
$ws.Range("N107").Value = -6197
$ws.Range("H137").Value = 5649.0713
$ws.Range("I137").Value = 4260.25
$ws.Range("J137").Value = 7500.8335
$ws.Range("K137").Value = 12780.75
$ws.Range("L137").Value = 22502.5005
$ws.Range("M137").Value = -10230.75
$ws.Range("N137").Value = -27602.5005
$ws.Range("H138").Value = 5751.092
$ws.Range("J138").Value = 6469.7188
$ws.Range("L138").Value = 19409.1564
$ws.Range("N138").Value = -29689.1564

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 24075
$ws.Range("I74").Value = 13751.667
$ws.Range("J74").Value = 106661.664
$ws.Range("K74").Value = 13751.667
$ws.Range("L74").Value = 106661.664
$ws.Range("M74").Value = -12877.667
$ws.Range("N74").Value = -108409.664
$ws.Range("H77").Value = 24075
$ws.Range("I77").Value = 13751.667
$ws.Range("J77").Value = 106661.664
$ws.Range("K77").Value = 68758.33499999999
$ws.Range("L77").Value = 533308.3200000001
$ws.Range("M77").Value = -64390.33499999999
$ws.Range("N77").Value = -542044.3200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H47").Value = 214679.5
$ws.Range("J47").Value = 214679.5
$ws.Range("L47").Value = 214679.5
$ws.Range("N47").Value = -215719.5
$ws.Range("H86").Value = 2127951.5
$ws.Range("I86").Value = 3402398.8
$ws.Range("K86").Value = 3402398.8
$ws.Range("M86").Value = -3401275.8
$ws.Range("H89").Value = 2127951.5
$ws.Range("I89").Value = 3402398.8
$ws.Range("K89").Value = 17011994
$ws.Range("M89").Value = -17006378

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 912810.4399999999
$ws.Range("I31").Value = 1002091.5
$ws.Range("K31").Value = 1002091.5
$ws.Range("M31").Value = -1001796.5
$ws.Range("H34").Value = 912810.4399999999
$ws.Range("I34").Value = 1002091.5
$ws.Range("K34").Value = 1002091.5
$ws.Range("M34").Value = -1001889.5
$ws.Range("H60").Value = 60057.145
$ws.Range("I60").Value = 35500
$ws.Range("J60").Value = 69880
$ws.Range("K60").Value = 35500
$ws.Range("L60").Value = 69880
$ws.Range("M60").Value = -34989
$ws.Range("N60").Value = -70902

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3656.4375
$ws.Range("I2").Value = 8007
$ws.Range("J2").Value = 272.66666
$ws.Range("K2").Value = 48042
$ws.Range("L2").Value = 1635.99996
$ws.Range("M2").Value = -47929
$ws.Range("N2").Value = -1861.99996
$ws.Range("H5").Value = 163519.4
$ws.Range("I5").Value = 800000
$ws.Range("J5").Value = 4399.25
$ws.Range("K5").Value = 2400000
$ws.Range("L5").Value = 13197.75
$ws.Range("M5").Value = -2399888
$ws.Range("N5").Value = -13421.75
$ws.Range("H92").Value = 667347.5600000001
$ws.Range("J92").Value = 1287
$ws.Range("L92").Value = 3861
$ws.Range("N92").Value = -6357
$ws.Range("H97").Value = 987.1667
$ws.Range("J97").Value = 605.5
$ws.Range("L97").Value = 1816.5
$ws.Range("N97").Value = -2808.5
$ws.Range("H135").Value = 163519.4
$ws.Range("I135").Value = 800000
$ws.Range("J135").Value = 4399.25
$ws.Range("K135").Value = 7200000
$ws.Range("L135").Value = 39593.25
$ws.Range("M135").Value = -7197465
$ws.Range("N135").Value = -44663.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1115341.1
$ws.Range("I80").Value = 838110.5
$ws.Range("J80").Value = 1669802.5
$ws.Range("K80").Value = 838110.5
$ws.Range("L80").Value = 1669802.5
$ws.Range("M80").Value = -837112.5
$ws.Range("N80").Value = -1671798.5
$ws.Range("H83").Value = 1115341.1
$ws.Range("I83").Value = 838110.5
$ws.Range("J83").Value = 1669802.5
$ws.Range("K83").Value = 4190552.5
$ws.Range("L83").Value = 8349012.5
$ws.Range("M83").Value = -4185560.5
$ws.Range("N83").Value = -8358996.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1198.3334
$ws.Range("I16").Value = 1197.25
$ws.Range("K16").Value = 1197.25
$ws.Range("M16").Value = -1027.25
$ws.Range("H61").Value = 3514.0244
$ws.Range("I61").Value = 1842.0952
$ws.Range("K61").Value = 1842.0952
$ws.Range("M61").Value = -1640.0952
$ws.Range("H113").Value = 3514.0244
$ws.Range("I113").Value = 1842.0952
$ws.Range("K113").Value = 1842.0952
$ws.Range("M113").Value = 327.9048

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H113").Value = 2748.3
$ws.Range("I113").Value = 2560.375
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 7681.125
$ws.Range("L113").Value = 10500
$ws.Range("M113").Value = -5511.125
$ws.Range("N113").Value = -14840
$ws.Range("H122").Value = 90917224
$ws.Range("I122").Value = 500001760
$ws.Range("J122").Value = 9556.111000000001
$ws.Range("K122").Value = 1500005280
$ws.Range("L122").Value = 28668.333
$ws.Range("M122").Value = -1500002830
$ws.Range("N122").Value = -33568.333
$ws.Range("H126").Value = 2080.6667
$ws.Range("I126").Value = 896.8
$ws.Range("K126").Value = 2690.4
$ws.Range("M126").Value = -220.3999999999996
$ws.Range("H132").Value = 20087.422
$ws.Range("I132").Value = 1763.7561
$ws.Range("J132").Value = 67041.81
$ws.Range("K132").Value = 5291.2683
$ws.Range("L132").Value = 201125.43
$ws.Range("M132").Value = -2761.2683
$ws.Range("N132").Value = -206185.43
